# Language sheet: add "Filipino" as a new row between Korean (row 5) and
# Hindi (row 7), i.e. fill in the previously-unused row 6 while leaving the
# rows below (7, 8) untouched, then move the active selection to A6.
#
# A plain Value write on A6 would leave the new cell with no cell style
# (s attribute omitted), but the target workbook has A6 using the same
# style index as the rows above it (s="1", the default-font style used
# throughout column A). Inserting a row copies the formatting of the row
# above into the new row, so we insert at row 6 (shifting 7/8 down to 8/9),
# write the value, and then delete the row that is now in position 7
# (the duplicate of the old row 7) to shift everything back up — restoring
# rows 7/8 to their original positions/content while keeping the inherited
# style on the new A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(6).Insert(-4121) | Out-Null      # xlShiftDown
$ws.Range("A6").Value = "Filipino"
$ws.Rows.Item(7).Delete(-4121) | Out-Null      # xlShiftUp, removes the duplicated old row 7

$ws.Range("A6").Select() | Out-Null
